# Appends the next day's (2026-01-15 / serial 46037) readings for both
# charging stations to the daily data sheet, then moves the selection the
# way the author left it (Excel re-persists the scroll/selection state of
# the active sheet on every save).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 四方坪站 (station 1) for 2026-01-15
$ws.Range("A30").Value = 46037
$ws.Range("B30").Value = "四方坪站"
$ws.Range("C30").Value = 14456.67
$ws.Range("D30").Value = 8909.13
$ws.Range("E30").Value = 3952.04
$ws.Range("F30").Value = 640

# 高岭站 (station 2) for 2026-01-15
$ws.Range("A31").Value = 46037
$ws.Range("B31").Value = "高岭站"
$ws.Range("C31").Value = 4160.84
$ws.Range("D31").Value = 3434.32
$ws.Range("E31").Value = 1069.49
$ws.Range("F31").Value = 154

# Scroll/selection bookkeeping that Excel writes into the sheetView on save.
$ws.Application.ActiveWindow.ScrollRow = 20
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("I30").Select()
